# Generate Report for Handback
# Regenerates the handback-status report with a new file pair:
#   984841b0-9947-491f-af4d-723d15d350d5.md -> b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md
#   a75b830c-eff0-4b78-b320-db89d208270a.md -> ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md
# plus refreshed timestamps / xliff correspondence file names.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md"
$ws1.Range("B2").Value = "e2e\b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md"
$ws1.Range("G2").Value = "2016-08-19 01:00:17"

$ws1.Range("A3").Value = "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md"
$ws1.Range("B3").Value = "e2e\ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md"
$ws1.Range("G3").Value = "2016-08-19 01:00:17"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e952eb5f123fab05b28e28113d26ff5ae272b6df/e2e/984841b0-9947-491f-af4d-723d15d350d5.md", "", "", "e2e\b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e952eb5f123fab05b28e28113d26ff5ae272b6df/e2e/a75b830c-eff0-4b78-b320-db89d208270a.md", "", "", "e2e\ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md"
$ws2.Range("G2").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-19 01:00:03"
$ws2.Range("I2").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md"
$ws2.Range("J2").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-19 01:00:31"

$ws2.Range("A3").Value = "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md"
$ws2.Range("G3").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-19 01:00:03"
$ws2.Range("I3").Value = "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md"
$ws2.Range("J3").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-19 01:00:31"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e952eb5f123fab05b28e28113d26ff5ae272b6df/e2e/984841b0-9947-491f-af4d-723d15d350d5.md", "", "", "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3a87e630be376dc609d69a6c3860243ee8f897b1/e2e/984841b0-9947-491f-af4d-723d15d350d5.md", "", "", "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e952eb5f123fab05b28e28113d26ff5ae272b6df/e2e/a75b830c-eff0-4b78-b320-db89d208270a.md", "", "", "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3a87e630be376dc609d69a6c3860243ee8f897b1/e2e/a75b830c-eff0-4b78-b320-db89d208270a.md", "", "", "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md"
$ws3.Range("G2").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-19 01:00:17"
$ws3.Range("I2").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md"
$ws3.Range("J2").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-19 01:00:39"

$ws3.Range("A3").Value = "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md"
$ws3.Range("G3").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-19 01:00:17"
$ws3.Range("I3").Value = "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md"
$ws3.Range("J3").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-19 01:00:39"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e952eb5f123fab05b28e28113d26ff5ae272b6df/e2e/984841b0-9947-491f-af4d-723d15d350d5.md", "", "", "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3fb76ab222c7b22194f56c45720dbcd8037cba10/e2e/984841b0-9947-491f-af4d-723d15d350d5.md", "", "", "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e952eb5f123fab05b28e28113d26ff5ae272b6df/e2e/a75b830c-eff0-4b78-b320-db89d208270a.md", "", "", "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3fb76ab222c7b22194f56c45720dbcd8037cba10/e2e/a75b830c-eff0-4b78-b320-db89d208270a.md", "", "", "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md")

Write-Host "Handback report regenerated."
